$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell="D2"; Value="332.85"},
    @{Cell="E2"; Value="1.12%"},
    @{Cell="G2"; Value="8"},
    @{Cell="D3"; Value="43.82"},
    @{Cell="E3"; Value="5.37%"},
    @{Cell="G3"; Value="8"},
    @{Cell="D4"; Value="5.848"},
    @{Cell="E4"; Value="3.86%"},
    @{Cell="G4"; Value="8"},
    @{Cell="D5"; Value="0.08344"},
    @{Cell="E5"; Value="1.99%"},
    @{Cell="G5"; Value="8"},
    @{Cell="D6"; Value="8.808"},
    @{Cell="E6"; Value="0.75%"},
    @{Cell="G6"; Value="8"},
    @{Cell="D7"; Value="1.989"},
    @{Cell="E7"; Value="-2.18%"},
    @{Cell="G7"; Value="8"},
    @{Cell="E8"; Value="-1.72%"},
    @{Cell="G8"; Value="8"},
    @{Cell="D9"; Value="0.9357"},
    @{Cell="E9"; Value="1.58%"},
    @{Cell="G9"; Value="8"},
    @{Cell="D10"; Value="0.1242"},
    @{Cell="E10"; Value="-2.65%"},
    @{Cell="G10"; Value="8"},
    @{Cell="D11"; Value="0.1950"},
    @{Cell="E11"; Value="-0.70%"},
    @{Cell="G11"; Value="8"},
    @{Cell="D12"; Value="0.09648"},
    @{Cell="E12"; Value="2.48%"},
    @{Cell="G12"; Value="8"},
    @{Cell="D13"; Value="0.04258"},
    @{Cell="E13"; Value="12.86%"},
    @{Cell="G13"; Value="8"},
    @{Cell="E14"; Value="0.67%"},
    @{Cell="G14"; Value="8"},
    @{Cell="D15"; Value="0.001306"},
    @{Cell="E15"; Value="0.33%"},
    @{Cell="G15"; Value="8"},
    @{Cell="D16"; Value="0.006017"},
    @{Cell="E16"; Value="-2.38%"},
    @{Cell="G16"; Value="8"},
    @{Cell="D17"; Value="3.495"},
    @{Cell="E17"; Value="1.48%"},
    @{Cell="G17"; Value="8"},
    @{Cell="D18"; Value="4.503"},
    @{Cell="E18"; Value="0.12%"},
    @{Cell="G18"; Value="8"},
    @{Cell="G19"; Value="8"},
    @{Cell="D20"; Value="8.826"},
    @{Cell="E20"; Value="6.23%"},
    @{Cell="G20"; Value="8"},
    @{Cell="D21"; Value="0.1362"},
    @{Cell="E21"; Value="-0.85%"},
    @{Cell="G21"; Value="8"},
    @{Cell="D22"; Value="0.2632"},
    @{Cell="E22"; Value="9.20%"},
    @{Cell="G22"; Value="8"},
    @{Cell="D23"; Value="0.04402"},
    @{Cell="E23"; Value="0.12%"},
    @{Cell="G23"; Value="8"},
    @{Cell="E24"; Value="0.09%"},
    @{Cell="G24"; Value="8"},
    @{Cell="D25"; Value="0.004426"},
    @{Cell="E25"; Value="2.82%"},
    @{Cell="G25"; Value="8"},
    @{Cell="D26"; Value="0.0001192"},
    @{Cell="E26"; Value="-0.75%"},
    @{Cell="G26"; Value="8"},
    @{Cell="D27"; Value="0.0003995"},
    @{Cell="G27"; Value="8"},
    @{Cell="G28"; Value="8"},
    @{Cell="G29"; Value="8"},
    @{Cell="G30"; Value="8"},
    @{Cell="G31"; Value="8"},
    @{Cell="G32"; Value="8"},
    @{Cell="G33"; Value="8"},
    @{Cell="G34"; Value="8"},
    @{Cell="G35"; Value="8"},
    @{Cell="G36"; Value="8"},
    @{Cell="G37"; Value="8"},
    @{Cell="G38"; Value="8"},
    @{Cell="D39"; Value="0.02799"},
    @{Cell="E39"; Value="1.07%"},
    @{Cell="G39"; Value="8"},
    @{Cell="D40"; Value="0.05782"},
    @{Cell="E40"; Value="7.00%"},
    @{Cell="G40"; Value="8"},
    @{Cell="D41"; Value="0.007916"},
    @{Cell="E41"; Value="-1.24%"},
    @{Cell="G41"; Value="8"},
    @{Cell="D42"; Value="0.1429"},
    @{Cell="E42"; Value="0.78%"},
    @{Cell="G42"; Value="8"},
    @{Cell="D43"; Value="0.009041"},
    @{Cell="E43"; Value="0.76%"},
    @{Cell="G43"; Value="8"},
    @{Cell="D44"; Value="0.002103"},
    @{Cell="E44"; Value="-3.14%"},
    @{Cell="G44"; Value="8"},
    @{Cell="D45"; Value="0.009871"},
    @{Cell="E45"; Value="-14.34%"},
    @{Cell="G45"; Value="8"},
    @{Cell="D46"; Value="0.00007251"},
    @{Cell="E46"; Value="9.49%"},
    @{Cell="G46"; Value="8"},
    @{Cell="E47"; Value="0.05%"},
    @{Cell="G47"; Value="8"},
    @{Cell="D48"; Value="0.003248"},
    @{Cell="E48"; Value="1.74%"},
    @{Cell="G48"; Value="8"},
    @{Cell="E49"; Value="-0.05%"},
    @{Cell="G49"; Value="8"},
    @{Cell="E50"; Value="0.05%"},
    @{Cell="G50"; Value="8"},
    @{Cell="E51"; Value="0.05%"},
    @{Cell="G51"; Value="8"}
)

foreach ($chg in $changes) {
    $rng = $ws.Range($chg.Cell)
    $oldStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $chg.Value
    $rng.Style = $oldStyle
}
